{"js": "// Apply the C3 conference justification-letter edits:\n//   1. \"April 20th\"  -> \"April 19th\"   (only the lone \"20\" run changes)\n//   2. \"CXSocial\"    -> \"Clarabridge Engage\"\n//   3. \"C3 Miami\"    -> \"C3 San Diego\"\n//   4. \"over 300\"    -> \"over 400\"\nconst body = context.document.body;\n\n// 1. The document contains five standalone \"20\" text runs:\n//      \"C3 20\" + \"20\" (-> \"C3 2020\"), \"April \" + \"20\" + \"th\" (-> \"April 20th\"),\n//      and \"... event of 20\" + \"20\" (-> \"...event of 2020\").\n//    Only the third occurrence (the one inside \"April 20th\") is changed by the\n//    edit, so target it specifically rather than doing a blanket replace that\n//    would also disturb the untouched \"2020\" occurrences / superscript \"th\" run.\nconst twenties = body.search(\"20\", { matchCase: true, matchWholeWord: false });\ntwenties.load(\"text\");\nawait context.sync();\ntwenties.items[2].insertText(\"19\", Word.InsertLocation.replace);\n\n// 2. \"CXSocial\" -> \"Clarabridge Engage\" (unique in the document).\nconst cxSocial = body.search(\"CXSocial\", { matchCase: true });\ncxSocial.load(\"text\");\nawait context.sync();\ncxSocial.items[0].insertText(\"Clarabridge Engage\", Word.InsertLocation.replace);\n\n// 3. \"Miami\" -> \"San Diego\" (unique in the document).\nconst miami = body.search(\"Miami\", { matchCase: true });\nmiami.load(\"text\");\nawait context.sync();\nmiami.items[0].insertText(\"San Diego\", Word.InsertLocation.replace);\n\n// 4. \"300\" -> \"400\" (unique in the document).\nconst threeHundred = body.search(\"300\", { matchCase: true });\nthreeHundred.load(\"text\");\nawait context.sync();\nthreeHundred.items[0].insertText(\"400\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Apply the C3 conference justification-letter edits:\n#   1. \"April 20th\"  -> \"April 19th\"   (only the lone \"20\" run changes)\n#   2. \"CXSocial\"    -> \"Clarabridge Engage\"\n#   3. \"C3 Miami\"    -> \"C3 San Diego\"\n#   4. \"over 300\"    -> \"over 400\"\n$d = $word.ActiveDocument\n\n# 1. The document contains five standalone \"20\" matches:\n#      \"C3 20\" + \"20\" (-> \"C3 2020\"), \"April \" + \"20\" + \"th\" (-> \"April 20th\"),\n#      and \"... event of 20\" + \"20\" (-> \"...event of 2020\").\n#    Only the third occurrence (the one inside \"April 20th\") is changed by the\n#    edit, so walk the matches and target that one specifically instead of doing\n#    a blanket replace that would also disturb the untouched \"2020\" occurrences.\n$count = 0\n$search = $d.Content\n$target = $null\nwhile ($search.Find.Execute(\"20\", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0)) {\n    $count = $count + 1\n    if ($count -eq 3) {\n        $target = $d.Range($search.Start, $search.End)\n        break\n    }\n    $search.Collapse(0)\n}\nif ($target -ne $null) {\n    $target.Text = \"19\"\n}\n\n# 2. \"CXSocial\" -> \"Clarabridge Engage\" (unique in the document).\n$d.Content.Find.Execute(\"CXSocial\", $false, $false, $false, $false, $false, $true, 1, $false, \"Clarabridge Engage\", 2) | Out-Null\n\n# 3. \"Miami\" -> \"San Diego\" (unique in the document).\n$d.Content.Find.Execute(\"Miami\", $false, $false, $false, $false, $false, $true, 1, $false, \"San Diego\", 2) | Out-Null\n\n# 4. \"300\" -> \"400\" (unique in the document).\n$d.Content.Find.Execute(\"300\", $false, $false, $false, $false, $false, $true, 1, $false, \"400\", 2) | Out-Null\n"}
